$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citi")

# Rows 3-5 hold values that look like scientific notation (e.g. "55E307570000"),
# so Excel would auto-parse them as numbers (overflowing to Infinity) unless the
# cells are pre-formatted as Text. Apply Text format first, assign the values,
# then switch to the target scientific display format (same format already used
# by column E on these rows), batching the format changes to avoid creating
# duplicate style entries.
$ws.Range("A3:A5").NumberFormat = "@"
$ws.Range("A3").Value = "55E307570000"
$ws.Range("A4").Value = "55E406695000"
$ws.Range("A5").Value = "55E397899000"
$ws.Range("A3:A5").NumberFormat = "0.00E+00"

# Row 6: C75XXX131 -> C75009131 (plain text, no style change)
$ws.Range("A6").Value = "C75009131"

# Row 7: C75XXX751 -> C75004751 (plain text, then apply scientific style like column E)
$ws.Range("A7").Value = "C75004751"
$ws.Range("A7").NumberFormat = "0.00E+00"

# Row 8: XXXX1802 (text) -> 38451802 (number)
$ws.Range("A8").Value = 38451802

# Row 9: clear A9 entirely (was XXXX1803)
$ws.Range("A9").ClearContents()

# Row 10: XXXX5752 (text) -> 10415752 (number)
$ws.Range("A10").Value = 10415752

# Row 11: XXXXXX9482 (text) -> 6866909482 (number)
$ws.Range("A11").Value = 6866909482

# Row 12: XXXXXX8419 (text) -> 9995698419 (number)
$ws.Range("A12").Value = 9995698419

# Update active cell selection to D15
$ws.Range("D15").Select()
